# This script applies the change described by the diff: splitting three
# run-level "wall of text" paragraphs (the Portuguese and English "Programa"
# paragraphs, and the "Bibliografia" paragraph) into multiple <w:t> pieces
# separated by manual line breaks (<w:br/>), matching the points between
# consecutive numbered items (Programa) or consecutive references
# (Bibliografia).
#
# Approach: use Range.Find.Execute with the "^l" special-character code in
# the replacement text, which Word interprets as a manual line break
# (w:br). Each Find/Replace call targets an exact, unique substring that
# straddles the desired break point, so the break lands in exactly the
# right place without disturbing the rest of the text.

$d = $word.ActiveDocument

function Insert-LineBreakAt($paragraphRange, [string]$findText, [string]$replaceText) {
    $ok = $paragraphRange.Find.Execute(
        $findText,    # FindText
        $false,       # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replaceText, # ReplaceWith
        1             # Replace (wdReplaceOne)
    )
    if (-not $ok) {
        throw "Could not find text: $findText"
    }
}

# Locate the three target paragraphs by their distinctive leading text,
# rather than hard-coded indices, so the script is resilient to any
# structural differences.
$paraPT = $null
$paraEN = $null
$paraBib = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("1. Qu") -and $t.Contains("amostragem ao tratamento de dados")) {
        $paraPT = $p
    } elseif ($t.StartsWith("1. Analytical chemistry in bioprocesses")) {
        $paraEN = $p
    } elseif ($t.StartsWith("- D. Harvey. Modern Analytical Chemistry")) {
        $paraBib = $p
    }
}

if ($null -eq $paraPT) { throw "Could not locate Portuguese Programa paragraph" }
if ($null -eq $paraEN) { throw "Could not locate English Programa paragraph" }
if ($null -eq $paraBib) { throw "Could not locate Bibliografia paragraph" }

# --- Portuguese "Programa" paragraph --------------------------------------
Insert-LineBreakAt $paraPT.Range "tratamento de dados.2. " "tratamento de dados.^l2. "
Insert-LineBreakAt $paraPT.Range "ação em fase sólida.3. " "ação em fase sólida.^l3. "
Insert-LineBreakAt $paraPT.Range "trometria de massas.4. " "trometria de massas.^l4. "

# --- English "Programa" paragraph ------------------------------------------
Insert-LineBreakAt $paraEN.Range "g to data treatment.2. " "g to data treatment.^l2. "
Insert-LineBreakAt $paraEN.Range "id phase extraction.3. " "id phase extraction.^l3. "
Insert-LineBreakAt $paraEN.Range "o mass spectrometry.4. " "o mass spectrometry.^l4. "

# --- "Bibliografia" paragraph ------------------------------------------
Insert-LineBreakAt $paraBib.Range "ill Companies, 2000.- A" "ill Companies, 2000.^l- A"
Insert-LineBreakAt $paraBib.Range "ollege London, 2004.- D" "ollege London, 2004.^l- D"
Insert-LineBreakAt $paraBib.Range "ning; 3ª Ed., 2023. - D" "ning; 3ª Ed., 2023. ^l- D"
Insert-LineBreakAt $paraBib.Range ". LTC; 9ª Ed., 2017.- D" ". LTC; 9ª Ed., 2017.^l- D"
Insert-LineBreakAt $paraBib.Range "rning, 2ª Ed., 2015.- O" "rning, 2ª Ed., 2015.^l- O"
Insert-LineBreakAt $paraBib.Range "nesp; 2ª Ed., 2011. - N" "nesp; 2ª Ed., 2011. ^l- N"
Insert-LineBreakAt $paraBib.Range " of Chemistry, 2018.- B" " of Chemistry, 2018.^l- B"

Write-Host "Done applying line breaks."
